# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) values on the last data row
# of the "zh-cn" and "de-de" sheets with fresh timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-03 07:02:00"
$wsZhCn.Range("G5").Value = "2016-03-03 07:02:49"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-03 07:02:12"
$wsDeDe.Range("G5").Value = "2016-03-03 07:03:07"
